$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z1").Value = [double]0.17
$ws.Range("Z2").Formula = "=0.17"
$ws.Range("Z3").Value = "0.17"
Write-Host $ws.Range("Z1").Value()
Write-Host $ws.Range("Z2").Value()
Write-Host $ws.Range("Z3").Value()
